$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Update column widths for columns E (5) and F (6)
# (ColumnWidth units get an implicit ~5/6 padding offset when stored as the
# OOXML "width" attribute, so back the values off by 5/6 to land on 17 / 18)
$ws.Columns.Item(5).ColumnWidth = 17 - 5/6
$ws.Columns.Item(6).ColumnWidth = 18 - 5/6

# Row 2 (OTROS)
$ws.Range("D2").Value = 68536.67999999999
$ws.Range("E2").Value = -68536.67999999999

# Row 3 (PORCELANATO)
$ws.Range("C3").Value = 15471.5593
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 15471.5593
$ws.Range("F3").Value = 0

# Row 4 (TOTAL)
$ws.Range("C4").Value = 15471.5593
$ws.Range("D4").Value = 68536.67999999999
$ws.Range("E4").Value = -53065.12069999999
$ws.Range("F4").Value = 4.429849549812344
